$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-23 10:02:38"
$wsZh.Range("H2").Value = "2016-03-23 10:03:20"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-23 10:02:46"
$wsDe.Range("H2").Value = "2016-03-23 10:03:33"
